# Updated cryptos list values (Price / Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.455.66"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  -0.76%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4601"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3822"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.64"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07911"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9704"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.04"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.825.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.883"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.039"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.95"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001027"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.458.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.347"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.048.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.062"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.301"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.97"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9539"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09295"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.244"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.313"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05932"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02191"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.043"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5785"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1838"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.252"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5478"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.865"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06636"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.53"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.039"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.90%  "

$ws.Range("E51").Value = "  -0.81%  "
